# Populate missing home/away score values that were left blank, and clean up
# a couple of cells that had picked up stray number formats (date format /
# a redundant "General" format) along the way - they should just use the
# same plain style as the rest of the score columns.

$wb = $excel.ActiveWorkbook

# ---- Division 3 --------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Division 3")
$ws3.Range("B2").Value = 1
$ws3.Range("C2").Value = 1

# ---- Division 4 --------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Division 4")
$ws4.Range("B2").Value = 1
$ws4.Range("C2").Value = 1
$ws4.Range("B3").Value = 1
$ws4.Range("C3").Value = 1

# ---- Cup -----------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Cup")

$ws5.Range("B2").Value = 1
$ws5.Range("C2").Value = 1

# B3/C3 previously carried a stray "General, applyNumberFormat" style - clear
# it first so the new value picks up the column's normal (unstyled) look.
$ws5.Range("B3:C3").Clear()
$ws5.Range("B3").Value = 1
$ws5.Range("C3").Value = 1

$ws5.Range("B4").Value = 1
$ws5.Range("C4").Value = 1

$ws5.Range("B5").Value = 1
$ws5.Range("C5").Value = 2

# B6 previously carried a stray date-number-format style - clear it first.
$ws5.Range("B6").Clear()
$ws5.Range("B6").Value = 1
$ws5.Range("C6").Value = 1

$ws5.Range("B7").Value = 3
$ws5.Range("C7").Value = 1

# B8 previously carried a stray date-number-format style - clear it first.
$ws5.Range("B8").Clear()
$ws5.Range("B8").Value = 1
$ws5.Range("C8").Value = 1

$ws5.Range("B9").Value = 1
$ws5.Range("C9").Value = 1

# B10 previously carried a stray date-number-format style - clear it first.
$ws5.Range("B10").Clear()
$ws5.Range("B10").Value = 5
$ws5.Range("C10").Value = 1

$ws5.Range("B11").Value = 1
$ws5.Range("C11").Value = 7

$ws5.Range("B12").Value = 1
$ws5.Range("C12").Value = 1

# B13 previously carried a stray date-number-format style - clear it first.
$ws5.Range("B13").Clear()
$ws5.Range("B13").Value = 1
$ws5.Range("C13").Value = 1

$ws5.Range("B14").Value = 1
$ws5.Range("C14").Value = 1

# ---- Selections / active sheet -------------------------------------------
# Division 3: cursor left on C40 after entering the scores.
$ws3.Range("C40").Select()

# Division 4: cursor left on C2.
$ws4.Range("C2").Select()

# Cup becomes the active tab, cursor parked on C8.
$ws5.Activate()
$ws5.Range("C8").Select()
